# Update "想去人数" (F column) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1361
$ws1.Range("F3").Value = 2886
$ws1.Range("F4").Value = 4
$ws1.Range("F5").Value = 265

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 5

# Sheet "全部类型" (All types) - aggregated view with same rows duplicated
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5
$ws4.Range("F3").Value = 1361
$ws4.Range("F4").Value = 2886
$ws4.Range("F5").Value = 4
$ws4.Range("F7").Value = 265
